$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All touched cells in columns B-E hold text (coin names, links, formatted
# price/volume strings) in the source data. Force text format on column D
# cells before assignment so Excel does not auto-convert numeric-looking
# strings (e.g. "585.34") into actual numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.912.14"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.288.72"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.34"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.53"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.427"
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.850.04"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.92"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.878.91"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.289.60"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.86"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.00"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.11"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +4.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.03"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.69"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  +3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.45"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.63"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.45"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "346.93"
$ws.Range("E44").Value = "  -5.10%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0690"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.628.81"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.78"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0284"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.08"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.33"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("E51").Value = "  -0.01%  "
